# Applies the "did some slight adjustments" commit:
#  - Adds a new row 13 on "Step 1 - Requirements check" with the text "["
#  - Highlights A15 (yellow) and A16 (green) on that sheet, matching the
#    "HIGHLIGHT IN GREEN WHAT APPLIES!" convention already used elsewhere
#  - Flips several checkbox (TRUE/FALSE) cells in column B to FALSE
#  - Switches the active sheet/selection from "Step 1" (A18) to
#    "Step 2 - Self assessment" (E9), which becomes the active tab

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# --- Step 1 sheet edits ---------------------------------------------------

# New row 13: a lone "[" character typed into A13 (column default style)
$ws1.Range("A13").Value = "["

# Highlight requirement bullets that apply (yellow / green), matching the
# "HIGHLIGHT IN GREEN WHAT APPLIES!" hint used elsewhere in the sheet
$ws1.Range("A15").Interior.Color = 65535
$ws1.Range("A16").Interior.Color = 5296274
$ws1.Range("A19").Interior.Color = 65535

# Uncheck several checkboxes (TRUE -> FALSE)
$ws1.Range("B16").Value = $false
$ws1.Range("B17").Value = $false
$ws1.Range("B23").Value = $false
$ws1.Range("B24").Value = $false

# --- View/selection state --------------------------------------------------

# Leave the selection on Step 1 at A13 (no longer the active/tabSelected sheet)
$ws1.Range("A13").Select()

# Step 2 becomes the active sheet with E9 selected
$ws2.Range("E9").Select()
$ws2.Activate()
